$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23 - entry #22
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "Monday, Jan 16"
$ws.Range("C23").Value = "9:05 AM"
$ws.Range("D23").Value = "W95177"
$ws.Range("E23").Value = "London"
$ws.Range("F23").Value = "(LTN)"
$ws.Range("G23").Value = "Wizz Air "
$ws.Range("H23").Value = "A320"
$ws.Range("I23").Value = "(G-WUKF)"
$ws.Range("J23").Value = "8:54 AM"
$ws.Range("K23").Borders.LineStyle = -4142
$ws.Range("L23").Value = "0 hours, -11 minutes"
$ws.Range("M23").Borders.LineStyle = -4142

# Row 24 - entry #23
$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "Monday, Jan 16"
$ws.Range("C24").Value = "10:25 AM"
$ws.Range("D24").Value = "FR2468"
$ws.Range("E24").Value = "London"
$ws.Range("F24").Value = "(STN)"
$ws.Range("G24").Value = "Ryanair "
$ws.Range("H24").Value = "B38M"
$ws.Range("I24").Value = "(EI-IFR)"
$ws.Range("J24").Value = "10:20 AM"
$ws.Range("K24").Borders.LineStyle = -4142
$ws.Range("L24").Value = "0 hours, -5 minutes"
$ws.Range("M24").Borders.LineStyle = -4142
